$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.268.14"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.891.36"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.61%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5071"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4048"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08311"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.23%  "
$ws.Range("D13").Value = "1.887.62"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.371"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.311"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001103"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06495"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.916"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "30.256.25"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.181"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("D26").Value = "2.100.10"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.264"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.106"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1041"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.010"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.725"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02447"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.326"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06443"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2151"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.186"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6404"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.578"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.219"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5964"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.131"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.641"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.214"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.24%  "
